# Updated capital structure database
# Apply new values to rows 2 and 3 (same new value for both rows in every
# changed column).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D"  = -0.0307
    "G"  = -2.172535211267606
    "H"  = -2.172535211267606
    "I"  = -2.27112676056338
    "J"  = -2.27112676056338
    "K"  = -7.14
    "L"  = -2.514084507042254
    "U"  = 0.047
    "V"  = 0.03012820512820513
    "W"  = 0.7547568710359407
    "X"  = 0.06185786519740737
    "Y"  = 0.6928990058385334
    "Z"  = -0.2985702270815812
    "AA" = 0.6780908326324643
    "AB" = 0.06185786519740737
    "AC" = 0.616232967435057
    "AG" = -0.047
    "AJ" = -0.03106411103767349
    "AK" = 0.002356244046723818
    "AL" = 3.27
    "AM" = 3.27
    "AO" = -1.972477064220183
    "AP" = 0.007355242566510173
    "AQ" = -1.972477064220183
}

foreach ($col in $updates.Keys) {
    $value = $updates[$col]
    $ws.Range("$col`2:$col`3").Value = $value
}
